# Actualización automática desde tarea programada
# Appends the latest sensor reading row coming from the scheduled task run,
# and refreshes the timestamp of the previous reading to its final captured
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: timestamp got its final precision once the task finished writing.
$ws.Range("A2").Value = 45866.04193254629

# Row 3: new sensor reading appended by the scheduled task.
$ws.Range("A3").Value = 45866.08356709986
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 31
$ws.Range("D3").Value = 13.78
$ws.Range("E3").Value = 91.63
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.06
$ws.Range("H3").Value = "ESE"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "02:00:20"
